$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: "populaire projecten toegevoegd" -> add new date/hour entries
$ws.Range("E18").Value = 41175
$ws.Range("F18").Value = 7

$ws.Range("E19").Value = 41193
$ws.Range("F19").Value = 2

$ws.Range("B24").Value = 41175
$ws.Range("C24").Value = 7

# Update the selection to match the recorded cursor position
$ws.Range("S21").Select()
